# Feat : FieldItem 오브젝트풀 적용
# Adds a new "gold" column (P) to the Entities sheet with per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Header for the new column (reuse the existing header style, same as O1)
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "gold"

# Per-row gold values (rows 2-19), default 5, bosses (rows 12 & 18) get 100
$goldValues = @{
    2  = 5
    3  = 5
    4  = 5
    5  = 5
    6  = 5
    7  = 5
    8  = 5
    9  = 5
    10 = 5
    11 = 5
    12 = 100
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 5
    18 = 100
    19 = 5
}

# Reuse the data-cell style (same as column O) for the new column's cells
$ws.Range("O2:O19").Copy()
$ws.Range("P2:P19").PasteSpecial(-4122)

foreach ($row in $goldValues.Keys) {
    $ws.Cells.Item($row, 16).Value = $goldValues[$row]
}

# Move the active selection to R17 (matches the recorded UI state in the diff)
$ws.Range("R17").Select()
